# Updates cryptocurrency Price (column D) and Volume(1h) (column E) data
# as refreshed by the scheduled GitHub Actions scraper run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.671.61"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").Value = "2.299.78"
$ws.Range("E3").Value = "  -0.50%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.88"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.96"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.06%  "
$ws.Range("E7").Value = "  -0.86%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.604"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.36%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.38"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.24%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0906"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.61%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.51"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.58%  "
$ws.Range("E13").Value = "  +1.24%  "
$ws.Range("E14").Value = "  +3.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.34"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.38%  "
$ws.Range("D16").Value = "2.648.79"
$ws.Range("E16").Value = "  -0.20%  "
$ws.Range("D17").Value = "2.299.48"
$ws.Range("E17").Value = "  +0.07%  "
$ws.Range("D18").Value = "42.607.83"
$ws.Range("E18").Value = "  +0.65%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.54"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.75%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.90"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +24.61%  "
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.98"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.40%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.53"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -3.82%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "266.51"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -5.40%  "
$ws.Range("E25").Value = "  -1.96%  "
$ws.Range("E26").Value = "  +0.47%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.91"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.27"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -5.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.61"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.64"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +11.89%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "37.15"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.26%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "165.40"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.15%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0881"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.31%  "
$ws.Range("E34").Value = "  -3.38%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.58"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -3.94%  "
$ws.Range("E36").Value = "  -2.66%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.54"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -3.15%  "
$ws.Range("E38").Value = "  +0.95%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.72"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.42%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.71"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.54%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.59"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +7.94%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "70.73"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.75%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "95.87"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -5.07%  "
$ws.Range("E44").Value = "  +0.48%  "
$ws.Range("E45").Value = "  +0.31%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.40"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.65%  "
$ws.Range("E47").Value = "  +3.75%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "79.97"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.65%  "
$ws.Range("D49").Value = "1.663.10"
$ws.Range("E49").Value = "  +3.46%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.27"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.75%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.84"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.00%  "
